$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 1
$ws.Range("E2").Value = 22
$ws.Range("G2").Value = "0.02122328812890018" -as [double]
$ws.Range("H2").Value = "0.09594376953620641" -as [double]
$ws.Range("I2").Value = "9.295927361208101e-10" -as [double]
$ws.Range("K2").Value = "5.223519391165876" -as [double]
$ws.Range("L2").Value = "[0.9180771413965978, 9.528961640935155]"
$ws.Range("M2").Value = "0.01777101687945337" -as [double]
$ws.Range("N2").Value = "0.01777101687945337" -as [double]
$ws.Range("O2").Value = "-1.559789745926464" -as [double]
$ws.Range("P2").Value = "[-2.7170531058073895, -0.4025263860455386]"
$ws.Range("Q2").Value = "0.008607579257601028" -as [double]
$ws.Range("R2").Value = "0.008607579257601028" -as [double]
$ws.Range("S2").Value = "12.10044248540324" -as [double]
$ws.Range("T2").Value = "[9.500468631848591, 14.70041633895788]"
$ws.Range("U2").Value = "4.440892098500626e-16" -as [double]
$ws.Range("V2").Value = "4.440892098500626e-16" -as [double]
$ws.Range("W2").Value = "5.461461461461461" -as [double]
$ws.Range("X2").Value = "1.409409409409407" -as [double]
$ws.Range("Y2").Value = "9.513513513513516" -as [double]

# Row 3 updates
$ws.Range("E3").Value = "22.68000000000011" -as [double]
$ws.Range("G3").Value = "0.00784715576651529" -as [double]
$ws.Range("H3").Value = "0.07875072427122859" -as [double]
$ws.Range("K3").Value = "4.178233688246032" -as [double]
$ws.Range("L3").Value = "[1.1169320894784374, 7.239535287013627]"
$ws.Range("M3").Value = "0.007644942468026716" -as [double]
$ws.Range("N3").Value = "0.01528988493605343" -as [double]
$ws.Range("O3").Value = "-1.58494764505431" -as [double]
$ws.Range("P3").Value = "[-2.6038425597320805, -0.5660527303765397]"
$ws.Range("Q3").Value = "0.002409833044479726" -as [double]
$ws.Range("R3").Value = "0.004819666088959451" -as [double]
$ws.Range("S3").Value = "12.8335951030887" -as [double]
$ws.Range("T3").Value = "[10.974823063598429, 14.692367142578973]"
$ws.Range("W3").Value = "5.72108108108111" -as [double]
$ws.Range("X3").Value = "2.043243243243256" -as [double]
$ws.Range("Y3").Value = "9.398918918918962" -as [double]
